# Append 15 new data rows (290-304) to the bottom of the sheet, matching
# the existing table's columns:
#   A Trigger_Level_High_Low | B Entry_Type | C Entry_Strike | D Strike_Type
#   E Expiry (date-formatted) | F Target | G Stop_Loss | H Qty | I Slicing
#   J Time_Interval | K Activation

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @(39000, "LIMIT", 39100, "PE", 45660, 39300, 38800, 4, 2, 5, 0),
    @(39000, "LIMIT", 39100, "PE", 45660, 39300, 38800, 4, 2, 5, 0),
    @(39000, "LIMIT", 39100, "PE", 45660, 39300, 38800, 4, 2, 5, 0),
    @(39300, "LIMIT", 39100, "CE", 45660, 39300, 38800, 4, 2, 5, 0),
    @(39300, "LIMIT", 39100, "CE", 45660, 38800, 39300, 4, 2, 5, 0),
    @(38750, "LIMIT", 38770, "PE", 45660, 38800, 38720, 4, 2, 5, 0),
    @(38750, "LIMIT", 38820, "PE", 45660, 38850, 38790, 4, 2, 5, 0),
    @(38850, "LIMIT", 38835, "CE", 45660, 38790, 38850, 4, 2, 5, 0),
    @(38850, "LIMIT", 38835, "CE", 45660, 38800, 38835, 4, 2, 5, 0),
    @(38850, "LIMIT", 38835, "CE", 45660, 38800, 38835, 4, 2, 5, 0),
    @(38850, "MARKET", 38845, "CE", 45660, 38835, 38855, 4, 2, 5, 0),
    @(39105, "LIMIT", 39120, "PE", 45660, 39200, 38900, 4, 2, 5, 0),
    @(39135, "LIMIT", 39120, "CE", 45660, 38800, 39150, 4, 2, 5, 0),
    @(39135, "LIMIT", 39120, "CE", 45660, 38800, 39150, 4, 2, 5, 0),
    @(39135, "LIMIT", 39120, "CE", 45660, 38800, 39150, 4, 2, 5, -1)
)

$startRow = 290
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $rowData = $rows[$i]
    for ($c = 0; $c -lt $rowData.Count; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $rowData[$c]
    }
    # Column E (5th column) carries the same date/time number format used
    # by the rest of the Expiry column.
    $ws.Cells.Item($r, 5).NumberFormat = "YYYY-MM-DD HH:MM:SS"
}
